# Fill in the "ACTUAL" working-hours value for the "OracleXE VM" task on
# both the Janeczek and Mair sheets, then move the active selection /
# active sheet to match the author's final cursor position.

$wb = $excel.ActiveWorkbook

$wsJaneczek = $wb.Worksheets.Item("Janeczek")
$wsMair     = $wb.Worksheets.Item("Mair")

# Enter the actual working hours (time-of-day fraction) for the
# "OracleXE VM" row on each sheet.
$wsJaneczek.Range("F8").Value = 0.086956018518518516
$wsMair.Range("F7").Value = 0.086956018518518516

# Update the selection on each sheet to reflect where the cursor ended up.
$wsJaneczek.Range("F16").Select()
$wsMair.Range("F7").Select()

# Make "Janeczek" the active (selected) sheet/tab.
$wsJaneczek.Activate()
